$d = $word.ActiveDocument

# 1. Title (Heading1) and the bold duplicate near the bottom - both occurrences share the
#    same old/new text, so Replace:=2 (wdReplaceAll) handles both.
$d.Content.Find.Execute(
    "Play Neptune's Champions Free Slot | Read Review", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Neptune’s Champions Slot Game Free", 2)

# 2. "What we like" bullets
$d.Content.Find.Execute(
    "Stunning graphics and immersive theme", $true, $false, $false, $false, $false,
    $true, 1, $false, "Medium-volatility for a balanced gaming experience", 2)

$d.Content.Find.Execute(
    "Free spins feature with increasing multipliers", $true, $false, $false, $false, $false,
    $true, 1, $false, "Visually stunning graphics and immersive theme", 2)

$d.Content.Find.Execute(
    "Medium-volatility for a balanced gameplay", $true, $false, $false, $false, $false,
    $true, 1, $false, "Potential for higher payouts with specific symbol combinations", 2)

$d.Content.Find.Execute(
    "Maximum payout potential up to 10,000x your bet", $true, $false, $false, $false, $false,
    $true, 1, $false, "Exciting free spins feature with increasing multipliers", 2)

# 3. "What we don't like" bullet
$d.Content.Find.Execute(
    "Higher payouts rely on specific symbol combinations", $true, $false, $false, $false, $false,
    $true, 1, $false, "No progressive jackpot feature", 2)

# 4. Italic meta description text
$d.Content.Find.Execute(
    "Read our review of Neptune's Champions slot game and play for free. Discover the stunning visuals, free spins feature, and 10,000x potential payout.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Neptune’s Champions slot game and play for free. Discover its features, volatility, and winning potential.", 2)
